$wb = $excel.ActiveWorkbook

# --- Worksheets (by name, to be robust against ordering) ---
$wsHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")   # sheet1
$wsHamba   = $wb.Worksheets.Item("Sedan_Hamba")     # sheet2
$wsMakhulu = $wb.Worksheets.Item("Bus_Makhulu")     # sheet3
$wsNone    = $wb.Worksheets.Item("None")            # sheet4

# --- Update fActuatorCutoff formulas: 1/0.2 -> 1/0.025 on the three brake sheets ---
$wsHambaLG.Range("H13").Formula = "=1/0.025"
$wsHambaLG.Range("H22").Formula = "=1/0.025"

$wsHamba.Range("H13").Formula = "=1/0.025"
$wsHamba.Range("H22").Formula = "=1/0.025"

$wsMakhulu.Range("H13").Formula = "=1/0.025"
$wsMakhulu.Range("H22").Formula = "=1/0.025"

# --- Update the per-sheet selected cell (bottom-right pane) ---
$wsHambaLG.Range("H22").Select()
$wsHamba.Range("H22").Select()
$wsMakhulu.Range("J22").Select()

# --- Switch the active/visible tab from "None" to "Bus_Makhulu" ---
$wsMakhulu.Activate()
